$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("A1").Value = "TEST"
